$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 17,3
$data[0,0] = 'Derrick White'
$data[0,1] = 'PG,SG'
$data[0,2] = 'Boston Celtics'
$data[1,0] = 'Bilal Coulibaly'
$data[1,1] = 'SG,SF'
$data[1,2] = 'Washington Wizards'
$data[2,0] = 'Kentavious Caldwell-Pope'
$data[2,1] = 'SG,SF'
$data[2,2] = 'Orlando Magic'
$data[3,0] = 'Herbert Jones'
$data[3,1] = 'SF,PF'
$data[3,2] = 'New Orleans Pelicans'
$data[4,0] = 'Julius Randle'
$data[4,1] = 'PF'
$data[4,2] = 'Minnesota Timberwolves'
$data[5,0] = 'Anthony Davis'
$data[5,1] = 'PF,C'
$data[5,2] = 'Los Angeles Lakers'
$data[6,0] = 'Damian Lillard'
$data[6,1] = 'PG'
$data[6,2] = 'Milwaukee Bucks'
$data[7,0] = 'Isaiah Hartenstein'
$data[7,1] = 'C'
$data[7,2] = 'Oklahoma City Thunder'
$data[8,0] = 'Brandon Miller'
$data[8,1] = 'SG,SF'
$data[8,2] = 'Charlotte Hornets'
$data[9,0] = 'Cameron Johnson'
$data[9,1] = 'SF,PF'
$data[9,2] = 'Brooklyn Nets'
$data[10,0] = 'Cade Cunningham'
$data[10,1] = 'PG,SG'
$data[10,2] = 'Detroit Pistons'
$data[11,0] = 'Malik Monk'
$data[11,1] = 'SG,SF'
$data[11,2] = 'Sacramento Kings'
$data[12,0] = 'Duncan Robinson'
$data[12,1] = 'SG,SF'
$data[12,2] = 'Miami Heat'
$data[13,0] = 'Kelly Olynyk'
$data[13,1] = 'C'
$data[13,2] = 'Toronto Raptors'
$data[14,0] = 'Bam Adebayo'
$data[14,1] = 'C'
$data[14,2] = 'Miami Heat'
$data[15,0] = 'Brandon Ingram'
$data[15,1] = 'SG,SF,PF'
$data[15,2] = 'New Orleans Pelicans'
$data[16,0] = 'LaMelo Ball'
$data[16,1] = 'PG,SG'
$data[16,2] = 'Charlotte Hornets'

$ws.Range("A2:C18").Value = $data
